$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of new row number -> old row number (where the data for the new row comes from)
$rowMap = @{
    2 = 2
    3 = 3
    4 = 4
    5 = 6
    6 = 5
    7 = 7
    8 = 8
    9 = 9
    10 = 10
    11 = 12
    12 = 11
    13 = 14
    14 = 13
    15 = 16
    16 = 15
    17 = 17
    18 = 18
    19 = 19
    20 = 20
    21 = 22
    22 = 23
    23 = 21
    24 = 24
    25 = 25
    26 = 29
    27 = 26
    28 = 27
    29 = 28
    30 = 30
    31 = 32
    32 = 35
    33 = 31
    34 = 34
    35 = 36
    36 = 37
    37 = 33
}

$firstRow = 2
$lastRow = 37

# Step 1: snapshot the "before" state of every row (values for A:R, formulas for S:Z)
$valuesByRow = @{}
$formulasByRow = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $valuesByRow[$r] = $ws.Range("A$r`:R$r").Value2
    $formulasByRow[$r] = $ws.Range("S$r`:Z$r").Formula
}

# Step 2: write back the rows in their new order, using the snapshotted data
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]

    $vals = $valuesByRow[$oldRow]
    $forms = $formulasByRow[$oldRow]

    $ws.Range("A$newRow`:R$newRow").Value2 = $vals
    $ws.Range("S$newRow`:Z$newRow").Formula = $forms

    # "Förändrad" (column C) is bumped by one day for every row
    $ws.Range("C$newRow").Value2 = 46078
}

Write-Host "Row reorder complete"
